$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.382.23'
$ws.Range('E2').Value = '  +0.04%  '
$ws.Range('D3').Value = '1.571.47'
$ws.Range('E3').Value = '  +0.23%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('E5').Value = '  +0.01%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '290.85'
$ws.Range('E6').Value = '  +0.13%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3757'
$ws.Range('E7').Value = '  +2.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '49.99'
$ws.Range('E8').Value = '  +1.18%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3418'
$ws.Range('E9').Value = '  +0.59%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07639'
$ws.Range('E10').Value = '  +0.50%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.150'
$ws.Range('E11').Value = '  -1.85%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.17'
$ws.Range('E13').Value = '  -0.18%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.020'
$ws.Range('E14').Value = '  -0.47%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.936'
$ws.Range('E15').Value = '  +0.62%  '
$ws.Range('D16').Value = '1.569.85'
$ws.Range('E16').Value = '  -0.86%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001130'
$ws.Range('E17').Value = '  -0.67%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '89.88'
$ws.Range('E18').Value = '  +0.71%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06731'
$ws.Range('E19').Value = '  -0.72%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.76'
$ws.Range('E21').Value = '  +1.48%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.198'
$ws.Range('E22').Value = '  -0.62%  '
$ws.Range('E23').Value = '  -0.44%  '
$ws.Range('D24').Value = '22.389.93'
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.396'
$ws.Range('E25').Value = '  +0.53%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.671'
$ws.Range('E26').Value = '  -10.27%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.18'
$ws.Range('E27').Value = '  +1.46%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '147.17'
$ws.Range('E28').Value = '  +1.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.029'
$ws.Range('E29').Value = '  +1.24%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '126.09'
$ws.Range('E30').Value = '  +0.52%  '
$ws.Range('D31').Value = '1.744.95'
$ws.Range('E31').Value = '  -0.80%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.145'
$ws.Range('E32').Value = '  -1.73%  '
$ws.Range('E33').Value = '  +0.72%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9831'
$ws.Range('E34').Value = '  -5.37%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.882'
$ws.Range('E35').Value = '  -4.23%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.08486'
$ws.Range('E36').Value = '  +0.30%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02540'
$ws.Range('E37').Value = '  -0.16%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.372'
$ws.Range('E38').Value = '  +9.93%  '
$ws.Range('B39').Value = 'Algorand'
$ws.Range('C39').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2320'
$ws.Range('E39').Value = '  -0.55%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06555'
$ws.Range('E40').Value = '  +0.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.410'
$ws.Range('E41').Value = '  -2.40%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6393'
$ws.Range('E42').Value = '  +0.29%  '
$ws.Range('E43').Value = '  -3.38%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.08'
$ws.Range('E44').Value = '  -2.75%  '
$ws.Range('B45').Value = 'Frax'
$ws.Range('C45').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.001'
$ws.Range('E45').Value = '  +0.11%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.777'
$ws.Range('E46').Value = '  -0.06%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5966'
$ws.Range('E47').Value = '  -0.68%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.290'
$ws.Range('E48').Value = '  +1.73%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.084'
$ws.Range('E49').Value = '  -2.46%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '125.38'
$ws.Range('E50').Value = '  +1.24%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07327'
$ws.Range('E51').Value = '  +0.54%  '

Write-Host "Applied crypto list updates"
